$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.03
$ws.Range("O2").Value = 1.25
$ws.Range("K4").Value = 1.91
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 2.7
$ws.Range("R4").Value = 1.44
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("X4").Value = 6.5
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AJ4").Value = 21
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 41
$ws.Range("G5").Value = 2.88
$ws.Range("I5").Value = 2.35
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 7.8
$ws.Range("J6").Value = 1.72
$ws.Range("L6").Value = 6.8
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.25
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.93
$ws.Range("W6").Value = 6.8
$ws.Range("X6").Value = 5.9
$ws.Range("Z6").Value = 7.2
$ws.Range("AF6").Value = 55
$ws.Range("AH6").Value = 20
$ws.Range("AJ6").Value = 20
$ws.Range("AL6").Value = 65
$ws.Range("AO6").Value = 5.7
$ws.Range("AQ6").Value = 14.5
$ws.Range("AU6").Value = 7.9
$ws.Range("AW6").Value = 9
$ws.Range("AY6").Value = 37
$ws.Range("N7").Value = 6.65
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 4.5
$ws.Range("L8").Value = 4.75
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("U8").Value = 1.8
$ws.Range("V8").Value = 1.91
$ws.Range("X8").Value = 8.5
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 251
$ws.Range("AI8").Value = 23
$ws.Range("AK8").Value = 51
$ws.Range("AN8").Value = 3.75
$ws.Range("AO8").Value = 9.5
$ws.Range("AX8").Value = 23
$ws.Range("AZ8").Value = 81
$ws.Range("BA8").Value = 101
$ws.Range("G9").Value = 1.5
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2.5
$ws.Range("AQ9").Value = 21
$ws.Range("AS9").Value = 81
$ws.Range("AW9").Value = 7.5
$ws.Range("AZ9").Value = 81
$ws.Range("BD9").Value = 151
$ws.Range("G10").Value = 2.6
$ws.Range("I10").Value = 2.75
$ws.Range("J10").Value = 3.4
$ws.Range("L10").Value = 3.5
$ws.Range("W10").Value = 7.5
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 11
$ws.Range("Z10").Value = 26
$ws.Range("AH10").Value = 7.5
$ws.Range("AN10").Value = 4.5
$ws.Range("AP10").Value = 29
$ws.Range("AS10").Value = 251
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 2.38
$ws.Range("W11").Value = 5.5
$ws.Range("AU11").Value = 9.5
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.53
